$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newline = [char]10

$newValue = '<rpc-reply message-id="urn:uuid:7a247e50-6bea-4a9f-aea2-11306f46386d">' + $newline +
'  <data>' + $newline +
'    <network-instances>' + $newline +
'      <network-instance>' + $newline +
'        <name>Prueba_LxVPN</name>' + $newline +
'        <config>' + $newline +
'          <name>Prueba_LxVPN</name>' + $newline +
'          <type>oc-ni-types:L3VRF</type>' + $newline +
'        </config>' + $newline +
'        <interfaces>' + $newline +
'          <interface>' + $newline +
'            <id>GigabitEthernet0/3/2</id>' + $newline +
'            <config>' + $newline +
'              <id>GigabitEthernet0/3/2</id>' + $newline +
'              <interface>GigabitEthernet0/3/2</interface>' + $newline +
'              <subinterface>0</subinterface>' + $newline +
'            </config>' + $newline +
'          </interface>' + $newline +
'        </interfaces>' + $newline +
'        <protocols>' + $newline +
'          <protocol>' + $newline +
'            <identifier>oc-pol-types:OSPF</identifier>' + $newline +
'            <name>22</name>' + $newline +
'            <config>' + $newline +
'              <identifier>oc-pol-types:OSPF</identifier>' + $newline +
'              <name>22</name>' + $newline +
'            </config>' + $newline +
'            <ospfv2>' + $newline +
'              <global>' + $newline +
'                <config>' + $newline +
'                  <router-id>172.16.1.3</router-id>' + $newline +
'                </config>' + $newline +
'              </global>' + $newline +
'            </ospfv2>' + $newline +
'          </protocol>' + $newline +
'          <protocol>' + $newline +
'            <identifier>oc-pol-types:STATIC</identifier>' + $newline +
'            <name>default</name>' + $newline +
'            <config>' + $newline +
'              <identifier>oc-pol-types:STATIC</identifier>' + $newline +
'              <name>default</name>' + $newline +
'            </config>' + $newline +
'          </protocol>' + $newline +
'          <protocol>' + $newline +
'            <identifier>oc-pol-types:DIRECTLY_CONNECTED</identifier>' + $newline +
'            <name>default</name>' + $newline +
'            <config>' + $newline +
'              <identifier>oc-pol-types:DIRECTLY_CONNECTED</identifier>' + $newline +
'              <name>default</name>' + $newline +
'            </config>' + $newline +
'          </protocol>' + $newline +
'        </protocols>' + $newline +
'      </network-instance>' + $newline +
'    </network-instances>' + $newline +
'  </data>' + $newline +
'</rpc-reply>' + $newline

$ws.Range("F2").Value = $newValue
